# "Calculo de impuestos" workbook: add the title row.
# The sheet originally has no data at all. We add a title in A2:J2,
# merged into a single cell and center-aligned, containing the text
# "Calculo de impuestos para la zona rural de Choluteca".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titleRange = $ws.Range("A2:J2")
$titleRange.Value = "Calculo de impuestos para la zona rural de Choluteca"
$titleRange.HorizontalAlignment = -4108  # xlCenter
$titleRange.Merge()

# Reflect the selection left behind on the sheet (A2, with A2:J2 selected)
$titleRange.Select() | Out-Null
